$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value even when it looks like a number,
# preserving the default "Normal" style (no NumberFormat residue).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "41.194.37"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "2.174.41"
$ws.Range("E3").Value = "  -2.59%  "

$ws.Range("E4").Value = "  -0.21%  "

Set-TextValue $ws.Range("D5") "248.17"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("E6").Value = "  -2.97%  "

Set-TextValue $ws.Range("D7") "65.41"
$ws.Range("E7").Value = "  -9.08%  "

$ws.Range("E8").Value = "  -0.01%  "

Set-TextValue $ws.Range("D9") "0.566"
$ws.Range("E9").Value = "  -4.59%  "

Set-TextValue $ws.Range("D10") "58.75"
$ws.Range("E10").Value = "  +0.76%  "

Set-TextValue $ws.Range("D11") "0.0926"
$ws.Range("E11").Value = "  -5.05%  "

Set-TextValue $ws.Range("D12") "35.61"
$ws.Range("E12").Value = "  -13.33%  "

Set-TextValue $ws.Range("D13") "0.104"
$ws.Range("E13").Value = "  -1.37%  "

Set-TextValue $ws.Range("D14") "6.83"
$ws.Range("E14").Value = "  -6.37%  "

$ws.Range("D15").Value = "2.497.95"
$ws.Range("E15").Value = "  -2.75%  "

Set-TextValue $ws.Range("D16") "14.29"
$ws.Range("E16").Value = "  -5.02%  "

Set-TextValue $ws.Range("D17") "0.848"
$ws.Range("E17").Value = "  -2.31%  "

$ws.Range("D18").Value = "2.179.69"
$ws.Range("E18").Value = "  -2.28%  "

$ws.Range("D19").Value = "41.088.62"
$ws.Range("E19").Value = "  -1.98%  "

$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  -3.79%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D21") "71.53"
$ws.Range("E21").Value = "  -2.21%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "6.07"
$ws.Range("E22").Value = "  -2.70%  "

Set-TextValue $ws.Range("D23") "228.20"
$ws.Range("E23").Value = "  -3.63%  "

Set-TextValue $ws.Range("D24") "2.02"
$ws.Range("E24").Value = "  -6.45%  "

Set-TextValue $ws.Range("D25") "3.79"
$ws.Range("E25").Value = "  -6.56%  "

$ws.Range("E26").Value = "  +0.09%  "

Set-TextValue $ws.Range("D27") "11.24"
$ws.Range("E27").Value = "  +2.95%  "

$ws.Range("E28").Value = "  -6.23%  "

Set-TextValue $ws.Range("D29") "3.71"
$ws.Range("E29").Value = "  -5.88%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D30") "167.97"
$ws.Range("E30").Value = "  -2.33%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D31") "2.09"
$ws.Range("E31").Value = "  -4.93%  "

Set-TextValue $ws.Range("D32") "20.16"
$ws.Range("E32").Value = "  -3.67%  "

Set-TextValue $ws.Range("D33") "0.121"
$ws.Range("E33").Value = "  -1.49%  "

Set-TextValue $ws.Range("D34") "5.66"
$ws.Range("E34").Value = "  +0.15%  "

$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("E36").Value = "  -3.91%  "

Set-TextValue $ws.Range("D37") "4.53"
$ws.Range("E37").Value = "  -4.19%  "

$ws.Range("E38").Value = "  +0.79%  "

Set-TextValue $ws.Range("D39") "24.27"
$ws.Range("E39").Value = "  -7.22%  "

Set-TextValue $ws.Range("D40") "0.0304"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("E41").Value = "  -4.94%  "

$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D42") "5.28"
$ws.Range("E42").Value = "  +7.88%  "

$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D43") "5.45"
$ws.Range("E43").Value = "  -9.10%  "

Set-TextValue $ws.Range("D44") "60.35"
$ws.Range("E44").Value = "  -11.52%  "

Set-TextValue $ws.Range("D45") "11.19"
$ws.Range("E45").Value = "  -5.87%  "

Set-TextValue $ws.Range("D46") "8.51"
$ws.Range("E46").Value = "  -3.43%  "

$ws.Range("E47").Value = "  -9.04%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.0991"
$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D49") "1.00"
$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("E51").Value = "  -4.54%  "

